$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.5730337078651685
$ws1.Range("C2").Value = 0.5407098121085595
$ws1.Range("D2").Value = 0.9700374531835206
$ws1.Range("E2").Value = 0.6943699731903485
$ws1.Range("F2").Value = 0.8371040723981901
$ws1.Range("G2").Value = 0.9412915851272016
$ws1.Range("H2").Value = 0.7734678561909971
$ws1.Range("I2").Value = 518
$ws1.Range("J2").Value = 440
$ws1.Range("K2").Value = 94
$ws1.Range("L2").Value = 16

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.8545454545454545
$ws2.Range("C2").Value = 0.1760299625468165
$ws2.Range("D2").Value = 0.2919254658385093

$ws2.Range("B3").Value = 0.5407098121085595
$ws2.Range("C3").Value = 0.9700374531835206
$ws2.Range("D3").Value = 0.6943699731903485

$ws2.Range("B4").Value = 0.5730337078651685
$ws2.Range("C4").Value = 0.5730337078651685
$ws2.Range("D4").Value = 0.5730337078651685
$ws2.Range("E4").Value = 0.5730337078651685

$ws2.Range("B5").Value = 0.6976276333270071
$ws2.Range("C5").Value = 0.5730337078651686
$ws2.Range("D5").Value = 0.4931477195144289

$ws2.Range("B6").Value = 0.6976276333270071
$ws2.Range("C6").Value = 0.5730337078651685
$ws2.Range("D6").Value = 0.4931477195144289

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 94
$ws3.Range("C2").Value = 440
$ws3.Range("B3").Value = 16
$ws3.Range("C3").Value = 518
